$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(22, 1).Value = 2000
$ws.Cells.Item(22, 2).Value = "struggle"
$ws.Cells.Item(22, 3).Value = 0.9749262332916433
$ws.Cells.Item(22, 4).Value = 1.300361778587099
$ws.Cells.Item(22, 5).Value = -6.267426431179062
$ws.Cells.Item(22, 6).Value = -0.4335615932941437
$ws.Cells.Item(22, 7).Value = 0.1406517177820205
$ws.Cells.Item(22, 8).Value = -0.8185594081878662

$ws.Cells.Item(23, 1).Value = 2100
$ws.Cells.Item(23, 2).Value = "struggle"
$ws.Cells.Item(23, 3).Value = -3.844243764877326
$ws.Cells.Item(23, 4).Value = 1.833226948976521
$ws.Cells.Item(23, 5).Value = -1.409952521324157
$ws.Cells.Item(23, 6).Value = 0.09498954564332961
$ws.Cells.Item(23, 7).Value = -0.7519751191139221
$ws.Cells.Item(23, 8).Value = -0.1093448773026466

$ws.Cells.Item(24, 1).Value = 2200
$ws.Cells.Item(24, 2).Value = "struggle"
$ws.Cells.Item(24, 3).Value = 0.4797788858413697
$ws.Cells.Item(24, 4).Value = -0.523662269115448
$ws.Cells.Item(24, 5).Value = -1.702465817332268
$ws.Cells.Item(24, 6).Value = 0.1846340149641037
$ws.Cells.Item(24, 7).Value = -1.312596678733826
$ws.Cells.Item(24, 8).Value = 0.0687223374843597

$ws.Cells.Item(25, 1).Value = 2300
$ws.Cells.Item(25, 2).Value = "struggle"
$ws.Cells.Item(25, 3).Value = 1.155098915100098
$ws.Cells.Item(25, 4).Value = 1.092013478279114
$ws.Cells.Item(25, 5).Value = 1.727226853370667
$ws.Cells.Item(25, 6).Value = 0.6478226184844971
$ws.Cells.Item(25, 7).Value = -0.9091202020645142
$ws.Cells.Item(25, 8).Value = -0.1838704347610473

$ws.Cells.Item(26, 1).Value = 2400
$ws.Cells.Item(26, 2).Value = "struggle"
$ws.Cells.Item(26, 3).Value = -1.098365545272828
$ws.Cells.Item(26, 4).Value = -0.6193101108074199
$ws.Cells.Item(26, 5).Value = 0.1845241859555233
$ws.Cells.Item(26, 6).Value = -0.1064432710409164
$ws.Cells.Item(26, 7).Value = -0.09178250283002851
$ws.Cells.Item(26, 8).Value = 0.0652098655700683

$ws.Cells.Item(27, 1).Value = 2500
$ws.Cells.Item(27, 2).Value = "struggle"
$ws.Cells.Item(27, 3).Value = -0.8518145084381094
$ws.Cells.Item(27, 4).Value = -0.03355145454406605
$ws.Cells.Item(27, 5).Value = 0.7549576908350003
$ws.Cells.Item(27, 6).Value = -0.042302418500185
$ws.Cells.Item(27, 7).Value = 0.3572034537792206
$ws.Cells.Item(27, 8).Value = 0.1937969923019409

$ws.Cells.Item(28, 1).Value = 2600
$ws.Cells.Item(28, 2).Value = "struggle"
$ws.Cells.Item(28, 3).Value = 0.405293345451355
$ws.Cells.Item(28, 4).Value = 0.8384262472391129
$ws.Cells.Item(28, 5).Value = 0.3231545425951481
$ws.Cells.Item(28, 6).Value = -0.2768746614456177
$ws.Cells.Item(28, 7).Value = 0.2338086664676666
$ws.Cells.Item(28, 8).Value = -0.1817324161529541

$ws.Cells.Item(29, 1).Value = 2700
$ws.Cells.Item(29, 2).Value = "struggle"
$ws.Cells.Item(29, 3).Value = 0.2438197135925255
$ws.Cells.Item(29, 4).Value = 0.4860433936119046
$ws.Cells.Item(29, 5).Value = -0.09267929568886754
$ws.Cells.Item(29, 6).Value = 0.0734565481543541
$ws.Cells.Item(29, 7).Value = 0.1968513280153274
$ws.Cells.Item(29, 8).Value = 0.1055269688367843

$ws.Cells.Item(30, 1).Value = 2800
$ws.Cells.Item(30, 2).Value = "struggle"
$ws.Cells.Item(30, 3).Value = -0.07322704792022328
$ws.Cells.Item(30, 4).Value = 0.1344193816185026
$ws.Cells.Item(30, 5).Value = -0.148086081258953
$ws.Cells.Item(30, 6).Value = -0.1507309973239898
$ws.Cells.Item(30, 7).Value = -0.0175623763352632
$ws.Cells.Item(30, 8).Value = 0.08170322328805921

$ws.Cells.Item(31, 1).Value = 2900
$ws.Cells.Item(31, 2).Value = "struggle"
$ws.Cells.Item(31, 3).Value = 0.4853796958923352
$ws.Cells.Item(31, 4).Value = 0.7483796477317809
$ws.Cells.Item(31, 5).Value = -0.1751452423632148
$ws.Cells.Item(31, 6).Value = -0.1734857261180877
$ws.Cells.Item(31, 7).Value = 0.0183259565383195
$ws.Cells.Item(31, 8).Value = 0.0250454749912023

